$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New version column header (K1)
$ws.Range("K1").Value = "V1.02"

# Mirror column J into new column K for rows 2-12 (same pass/fail values)
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 11).Value = $ws.Cells.Item($r, 10).Value2
}

# New row 13 - AppearanceID test
$ws.Range("A13").Value = "AppearanceID"
$ws.Range("B13").Value = "Ensure only one sprite is set for each appearance ID`nTest can only occur at runtime due to Unity limitations, test will be disabled in build for optimisation"
$ws.Range("C13").Value = "-"
$ws.Range("D13").Value = "-"
$ws.Range("E13").Value = "Yes"
$ws.Range("F13").Value = "No"
$ws.Range("H13").Value = "-"
$ws.Range("I13").Value = "-"
$ws.Range("J13").Value = "-"
$ws.Range("K13").Value = "-"

$ws.Rows.Item(13).RowHeight = 75

# Extend conditional formatting range to include column K
$fc = $ws.Range("H1:J1048576").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("H1:K1048576"))

# Update selection to match target state
$ws.Range("L11").Select()
